$d = $word.ActiveDocument

# Replace "NAAM:" with "1" in the document body
$d.Content.Find.Execute("NAAM:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1", 2)

$d.Save()
